# Apply updated dSF (column F) values as part of "repull data, push all data,
# mean calculation" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new F column value
$updates = @{
    2  = -8
    3  = -1
    7  = -8
    8  = 5
    9  = -6
    10 = -2
    11 = -1
    12 = 1
    14 = -2
    15 = -4
    16 = -6
    17 = -2
    18 = -2
    19 = -8
    22 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
